# Add 2022-Q1 sheet (new quarterly fund-holding snapshot) and update the
# "总计" (totals) summary sheet with a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$new = $wb.Worksheets.Add($zj)
$new.Name = "2022-Q1"
# Re-resolve "总计" by name: the engine's worksheet handles are
# position-bound, and inserting a sheet in front of it shifts its index,
# so the original $zj reference would now (incorrectly) point at $new.
$zj = $wb.Worksheets.Item("总计")

# Match this workbook's page-margin convention (0.75/0.75/1/1/0.5/0.5 in,
# expressed in points) instead of Excel's generic new-sheet defaults.
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Data rows 2-10: column A (index) and column H (rank) are real numbers;
# columns D/E/F/G are text (leading apostrophe keeps "54.05" etc. from
# being coerced into numeric values); B/C are plain text.
$rows = @(
    @(0, "167301", "方正富邦中证保险主题指数（LOF）", "54.05", "93.05", "7.58", "4.0970", 5),
    @(1, "004702", "南方金融主题灵活配置混合", "30.93", "89.81", "4.40", "1.3609", 9),
    @(2, "163407", "兴全沪深300指数增强(LOF)A", "41.45", "95.41", "2.08", "0.8622", 10),
    @(3, "257040", "国联安红利混合", "1.10", "72.31", "4.48", "0.0493", 7),
    @(4, "007230", "兴全沪深300指数增强(LOF)C", "1.38", "95.41", "2.08", "0.0287", 10),
    @(5, "350001", "天治财富增长混合", "0.98", "69.00", "2.86", "0.0280", 4),
    @(6, "516720", "浦银安盛中证ESG 120策略交易型开放式指数证券投资基金", "0.86", "96.67", "2.48", "0.0213", 8),
    @(7, "012977", "瑞达鑫红量化6个月持有混合型证券投资基金A", "1.04", "94.56", "1.02", "0.0106", 5),
    @(8, "012978", "瑞达鑫红量化6个月持有混合型证券投资基金C", "0.17", "94.56", "1.02", "0.0017", 5)
)

$r = 2
foreach ($row in $rows) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = "'" + $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = "'" + $row[3]
    $new.Cells.Item($r, 5).Value = "'" + $row[4]
    $new.Cells.Item($r, 6).Value = "'" + $row[5]
    $new.Cells.Item($r, 7).Value = "'" + $row[6]
    $new.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Formatting: header row (B1:H1) and index column (A2:A10) use the bold +
# thin-border + centered style (style index 2 in the original workbook).
# Copy it from the "总计" sheet's own header/index cells so the same style
# entry is reused.
$zj.Range("B1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

$zj.Range("A2").Copy()
$new.Range("A2:A10").PasteSpecial(-4122)

$new.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row 2 for 2022-Q1, pushing the
#    existing quarters down.
# ---------------------------------------------------------------------
$zj.Rows.Item(2).Insert()

$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("B3:D3").Copy()
$zj.Range("B2:D2").PasteSpecial(-4122)

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 9
$zj.Cells.Item(2, 4).Value = 6.46

# Column A is a simple 0-based row counter, so every row below the new
# insertion point shifts up by one (0,1,2,3,4 -> 1,2,3,4,5).
$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(6, 1).Value = 4
$zj.Cells.Item(7, 1).Value = 5
